# Add two new locator rows ("modalContent" and "modalContentCloseBtn")
# into the digital_coupons_page sheet, just before the existing
# "loadedText" row (currently row 18), pushing the rows below it down
# by two.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("digital_coupons_page")

# Insert two blank rows at 18:19 - everything from the old row 18
# downward shifts to 20 downward.
$ws.Range("A18:A19").EntireRow.Insert()

# New row 18: modalContent locator (write value column first so the
# shared-string table order matches the source workbook).
$ws.Range("B18").Value = "//div[contains(@class,'modal-content')]"
$ws.Range("A18").Value = "modalContent"
$ws.Range("C18").Value = "xpath"

# New row 19: modalContentCloseBtn locator.
$ws.Range("B19").Value = "//div[contains(@class,'modal-content')]//button[contains(text(),'Close')]"
$ws.Range("A19").Value = "modalContentCloseBtn"
$ws.Range("C19").Value = "xpath"

# Match the saved view state: selection resting on F30 with no special
# top-left scroll position.
$ws.Activate()
$ws.Range("F30").Select()
